$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New result table (6 data rows) replacing the previous 2-row table,
# reflecting the updated Natmi analysis (adds "ECs" cluster; recomputes
# all sending/target cluster combinations for Bmp2-Bmpr1b).
$data = @(
    @("ECs","Bmp2","Bmpr1b","FAPs",2,0.6666666666666666,9.163165666666666,27.489497,0.2800251397703982,0.2800251397703982,3,1,1.292736333333333,3.878209,0.9024289107934966,0.9024289107934966,11.84555718565255,106.610014670873,0.2527027818777971,0.2527027818777971),
    @("ECs","Bmp2","Bmpr1b","sCs",2,0.6666666666666666,9.163165666666666,27.489497,0.2800251397703982,0.2800251397703982,2,0.6666666666666666,0.1397713333333333,0.419314,0.09757108920650338,0.09757108920650338,1.280747882784222,11.526730945058,0.0273223578926011,0.0273223578926011),
    @("FAPs","Bmp2","Bmpr1b","FAPs",3,1,12.06704733333333,36.201142,0.3687673822623249,0.3687673822623249,3,1,1.292736333333333,3.878209,0.9024289107934966,0.9024289107934966,15.59951052385311,140.395594714678,0.3327863471111588,0.3327863471111588),
    @("FAPs","Bmp2","Bmpr1b","sCs",3,1,12.06704733333333,36.201142,0.3687673822623249,0.3687673822623249,2,0.6666666666666666,0.1397713333333333,0.419314,0.09757108920650338,0.09757108920650338,1.686627295176444,15.179645656588,0.03598103515116603,0.03598103515116603),
    @("sCs","Bmp2","Bmpr1b","FAPs",3,1,11.49244066666667,34.477322,0.3512074779672769,0.351207477967277,3,1,1.292736333333333,3.878209,0.9024289107934966,0.9024289107934966,14.85669560847755,133.710260476298,0.3169397818045407,0.3169397818045407),
    @("sCs","Bmp2","Bmpr1b","sCs",3,1,11.49244066666667,34.477322,0.3512074779672769,0.351207477967277,2,0.6666666666666666,0.1397713333333333,0.419314,0.09757108920650338,0.09757108920650338,1.606313755234222,14.456823797108,0.03426769616273625,0.03426769616273626)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $row[$j]
    }
}
